$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 15 (the "checksum" sub-field row under documents[].file), which shifts
# all subsequent rows (16-77) up by one, so before-row-77 disappears entirely.
$ws.Rows.Item(15).Delete()
